# Scheduled-runner update: refresh market-board derived figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) on a handful of Leve
# rows across the ALC / ARM / BSM / CRP / CUL / GSM / LTW / WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1158.8334
$ws.Range("I28").Value = 1545
$ws.Range("J28").Value = 883
$ws.Range("K28").Value = 1545
$ws.Range("L28").Value = 883
$ws.Range("M28").Value = -1060
$ws.Range("N28").Value = -1853
$ws.Range("H53").Value = 215.5
$ws.Range("I53").Value = 131.875
$ws.Range("J53").Value = 550
$ws.Range("K53").Value = 131.875
$ws.Range("L53").Value = 550
$ws.Range("M53").Value = 505.125
$ws.Range("N53").Value = -1824
$ws.Range("H74").Value = 10597.143
$ws.Range("I74").Value = 9030
$ws.Range("K74").Value = 9030
$ws.Range("M74").Value = -8094
$ws.Range("H77").Value = 10597.143
$ws.Range("I77").Value = 9030
$ws.Range("K77").Value = 45150
$ws.Range("M77").Value = -40470
$ws.Range("H86").Value = 5730.8
$ws.Range("I86").Value = 4885
$ws.Range("J86").Value = 6999.5
$ws.Range("K86").Value = 4885
$ws.Range("L86").Value = 6999.5
$ws.Range("M86").Value = -3762
$ws.Range("N86").Value = -9245.5
$ws.Range("H89").Value = 5730.8
$ws.Range("I89").Value = 4885
$ws.Range("J89").Value = 6999.5
$ws.Range("K89").Value = 24425
$ws.Range("L89").Value = 34997.5
$ws.Range("M89").Value = -18809
$ws.Range("N89").Value = -46229.5
$ws.Range("H98").Value = 1998
$ws.Range("I98").Value = 1998
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1998
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -500
$ws.Range("N98").ClearContents()
$ws.Range("H111").Value = 2500
$ws.Range("I111").Value = 1500
$ws.Range("J111").Value = 3000
$ws.Range("K111").Value = 4500
$ws.Range("L111").Value = 9000
$ws.Range("M111").Value = -1433
$ws.Range("N111").Value = -15134
$ws.Range("H112").Value = 2131
$ws.Range("I112").Value = 1800
$ws.Range("J112").Value = 2351.6667
$ws.Range("K112").Value = 5400
$ws.Range("L112").Value = 7055.000100000001
$ws.Range("M112").Value = -4292
$ws.Range("N112").Value = -9271.000100000001
$ws.Range("H122").Value = 1998
$ws.Range("I122").Value = 1998
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5994
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3544
$ws.Range("N122").ClearContents()
$ws.Range("H138").Value = 2714.75
$ws.Range("J138").Value = 2986.3333
$ws.Range("L138").Value = 8958.999899999999
$ws.Range("N138").Value = -19238.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 12619.588
$ws.Range("J44").Value = 12619.588
$ws.Range("L44").Value = 12619.588
$ws.Range("N44").Value = -13595.588
$ws.Range("H55").Value = 52999.5
$ws.Range("I55").Value = 6000
$ws.Range("K55").Value = 6000
$ws.Range("M55").Value = -5685
$ws.Range("H122").Value = 1469.2759
$ws.Range("I122").Value = 1203.6522
$ws.Range("K122").Value = 3610.9566
$ws.Range("M122").Value = -1160.9566

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5203.15
$ws.Range("I86").Value = 3093.1
$ws.Range("K86").Value = 3093.1
$ws.Range("M86").Value = -1970.1
$ws.Range("H89").Value = 5203.15
$ws.Range("I89").Value = 3093.1
$ws.Range("K89").Value = 15465.5
$ws.Range("M89").Value = -9849.5
$ws.Range("H94").Value = 1480.0667
$ws.Range("I94").Value = 1492.4615
$ws.Range("J94").Value = 1399.5
$ws.Range("K94").Value = 1492.4615
$ws.Range("L94").Value = 1399.5
$ws.Range("M94").Value = -1041.4615
$ws.Range("N94").Value = -2301.5
$ws.Range("H105").Value = 1928.4
$ws.Range("I105").Value = 2048.2
$ws.Range("J105").Value = 1808.6
$ws.Range("K105").Value = 2048.2
$ws.Range("L105").Value = 1808.6
$ws.Range("M105").Value = -301.1999999999998
$ws.Range("N105").Value = -5302.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 429.1111
$ws.Range("J22").Value = 500.33334
$ws.Range("L22").Value = 500.33334
$ws.Range("N22").Value = -1200.33334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 168673.42
$ws.Range("J4").Value = 1761.25
$ws.Range("L4").Value = 5283.75
$ws.Range("N4").Value = -5507.75
$ws.Range("H18").Value = 2912.8572
$ws.Range("I18").Value = 130
$ws.Range("K18").Value = 390
$ws.Range("M18").Value = -221
$ws.Range("H60").Value = 1562.5
$ws.Range("I60").Value = 150.33333
$ws.Range("J60").Value = 2033.2222
$ws.Range("K60").Value = 450.99999
$ws.Range("L60").Value = 6099.6666
$ws.Range("M60").Value = -199.99999
$ws.Range("N60").Value = -6601.6666
$ws.Range("H61").Value = 133.77777
$ws.Range("I61").Value = 138.625
$ws.Range("K61").Value = 415.875
$ws.Range("M61").Value = -200.875
$ws.Range("H81").Value = 1511
$ws.Range("J81").Value = 1996.3334
$ws.Range("L81").Value = 5989.0002
$ws.Range("N81").Value = -8235.0002
$ws.Range("H84").Value = 1511
$ws.Range("J84").Value = 1996.3334
$ws.Range("L84").Value = 17967.0006
$ws.Range("N84").Value = -29199.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7798.4
$ws.Range("J80").Value = 5000
$ws.Range("L80").Value = 5000
$ws.Range("N80").Value = -6996
$ws.Range("H83").Value = 7798.4
$ws.Range("J83").Value = 5000
$ws.Range("L83").Value = 25000
$ws.Range("N83").Value = -34984
$ws.Range("H102").Value = 1472.1818
$ws.Range("I102").Value = 1472.1818
$ws.Range("K102").Value = 1472.1818
$ws.Range("M102").Value = 149.8181999999999
$ws.Range("H107").Value = 4714.143
$ws.Range("I107").Value = 3999.5
$ws.Range("K107").Value = 3999.5
$ws.Range("M107").Value = -2079.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 47933
$ws.Range("J98").Value = 47933
$ws.Range("L98").Value = 47933
$ws.Range("N98").Value = -53923
$ws.Range("H127").Value = 62305
$ws.Range("J127").Value = 62305
$ws.Range("L127").Value = 62305
$ws.Range("N127").Value = -72225

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9285.571
$ws.Range("I81").Value = 6625
$ws.Range("J81").Value = 12833
$ws.Range("K81").Value = 13250
$ws.Range("L81").Value = 25666
$ws.Range("M81").Value = -12189
$ws.Range("N81").Value = -27788
$ws.Range("H84").Value = 9285.571
$ws.Range("I84").Value = 6625
$ws.Range("J84").Value = 12833
$ws.Range("K84").Value = 66250
$ws.Range("L84").Value = 128330
$ws.Range("M84").Value = -60946
$ws.Range("N84").Value = -138938
$ws.Range("H96").Value = 1446.1428
$ws.Range("I96").Value = 1594.6
$ws.Range("K96").Value = 1594.6
$ws.Range("M96").Value = -221.5999999999999
$ws.Range("H100").Value = 928.0909
$ws.Range("I100").Value = 820.9
$ws.Range("K100").Value = 1641.8
$ws.Range("M100").Value = -1100.8
$ws.Range("H117").Value = 54999.668
$ws.Range("J117").Value = 54999.668
$ws.Range("L117").Value = 54999.668
$ws.Range("N117").Value = -64177.668
$ws.Range("H122").Value = 4250.357
$ws.Range("I122").Value = 4292.1665
$ws.Range("J122").Value = 3999.5
$ws.Range("K122").Value = 12876.4995
$ws.Range("L122").Value = 11998.5
$ws.Range("M122").Value = -10426.4995
$ws.Range("N122").Value = -16898.5
$ws.Range("H126").Value = 5168.952
$ws.Range("I126").Value = 2676.2
$ws.Range("K126").Value = 8028.599999999999
$ws.Range("M126").Value = -5558.599999999999
